$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.970.69"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.516.41"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'608.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'147.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("D7").Value = "3.515.76"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E11").Value = "  +4.73%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "4.110.09"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "'31.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "3.522.43"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "66.984.01"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'10.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.79%  "
$ws.Range("D20").Value = "'6.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'438.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").Value = "'79.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "3.658.05"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "'9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "'8.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'25.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "3.508.13"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "'5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'173.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -9.45%  "
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'46.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'28.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'7.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'2.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("E51").Value = "  +0.64%  "
